# RoboCUP SSL Hardware BOM — add two Break Beam sensor options (5mm + 3mm)
# to the "Break Beam BOM" sheet, let the "Complete BOM" rollup sheet pick
# the values up through its existing cross-sheet formulas, and drop the
# stale "Haven't received a quotation" comment that no longer applies.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Break Beam BOM" sheet: fill in the Secondary (D, 3mm LEDs) and
#    Primary (B, 5mm LEDs) component columns. The 3mm (D) column was
#    entered first in the authored workbook, then the 5mm (B) column —
#    kept in that order so shared strings line up the same way.
# ---------------------------------------------------------------------
$bbWs = $wb.Worksheets.Item("Break Beam BOM")

# --- Column D: IR Break Beam Sensors - 3mm LEDs ---
$bbWs.Range("D1").Value = "IR Break Beam Sensors with Premium Wire Header Ends - 3mm LEDs"

$d2Text = "IR Break Beam Sensors with Premium Wire Header Ends - 3mm LEDs : ID 2167 : Adafruit Industries, Unique & fun DIY electronics and kits"
$bbWs.Range("D2").Value = $d2Text
$bbWs.Hyperlinks.Add($bbWs.Range("D2"), "https://www.adafruit.com/product/2167", [Type]::Missing, [Type]::Missing, $d2Text) | Out-Null

$bbWs.Range("D3").Value = "(20*10*8)"
$bbWs.Range("D4").Value = 1
$bbWs.Range("D5").Formula = "=10.24*2.95"
# D6 already holds =D4*D5 from the template; recalculated automatically.

# --- Column B: IR Break Beam Sensor - 5mm LEDs ---
$bbWs.Range("B1").Value = "IR Break Beam Sensor with Premium Wire Header Ends - 5mm LEDs"

$b2Text = "IR Break Beam Sensor with Premium Wire Header Ends - 5mm LEDs : ID 2168 : Adafruit Industries, Unique & fun DIY electronics and kits"
$bbWs.Range("B2").Value = $b2Text
$bbWs.Hyperlinks.Add($bbWs.Range("B2"), "https://www.adafruit.com/product/2168", [Type]::Missing, [Type]::Missing, $b2Text) | Out-Null

$bbWs.Range("B3").Value = "(20*10*8)"
$bbWs.Range("B4").Value = 1
$bbWs.Range("B5").Formula = "=10.24*5.95"
$bbWs.Range("B6").Formula = "=B4*B5"

# Widen the two newly-populated columns so the long descriptions are
# readable (mirrors the "bestFit" auto-sizing applied in the authored
# workbook).
$bbWs.Columns.Item(2).ColumnWidth = 112.8
$bbWs.Columns.Item(4).ColumnWidth = 55.6

# Leave the cursor where the author left it after entering the data.
$bbWs.Activate()
$bbWs.Range("B16").Select()

# ---------------------------------------------------------------------
# 2. "Complete BOM" sheet: the Break Beam column (C) is pulled in purely
#    via formulas already in the sheet, so no manual writes are needed
#    there — just widen column C to fit the long strings now flowing
#    into it and restore the author's final selection / active sheet.
# ---------------------------------------------------------------------
$mainWs = $wb.Worksheets.Item("Complete BOM")
$mainWs.Columns.Item(3).ColumnWidth = 112.8

# ---------------------------------------------------------------------
# 3. Drop the obsolete "Haven't received a quotation" comment on D20 —
#    a price is now available, so the note no longer applies.
# ---------------------------------------------------------------------
if ($mainWs.Comments.Count -gt 0) {
    $mainWs.Comments.Item(1).Delete()
}

# Restore "Complete BOM" as the visible/active sheet with its last
# selection, matching the saved workbook state.
$mainWs.Activate()
$mainWs.Range("C36").Select()

# Make sure every cross-sheet formula (Complete BOM pulling from Break
# Beam BOM) is recalculated before the workbook is saved.
$wb.Application.Calculate()
